$wb = $excel.ActiveWorkbook

# --- Reorder sheets ---
# Before: ... GCost, RCost, Zone, AGCR
# After:  ... GCost, Zone, AGCR(->SFR), RCost(->SFRCost)
# Move RCost to the end (after AGCR) so the order becomes Zone, AGCR, RCost.
$rcost = $wb.Worksheets.Item("RCost")
$agcr = $wb.Worksheets.Item("AGCR")
$rcost.Move($null, $agcr)

# --- Rename sheets ---
# NOTE: worksheet handles obtained before a Move/reorder become stale
# (they resolve by position), so re-fetch by name after reordering.
$wb.Worksheets.Item("AGCR").Name = "SFR"
$wb.Worksheets.Item("RCost").Name = "SFRCost"

# --- Active sheet / selected tab ---
# The saved file shows the last sheet (SFRCost, formerly RCost) as the
# active/selected tab, with cell Q23 selected.
$sfrcost = $wb.Worksheets.Item("SFRCost")
$sfrcost.Activate() | Out-Null
$sfrcost.Range("Q23").Select() | Out-Null
